# Update the "handback-status" timestamps that were regenerated when the
# handback report was (re-)generated. This mirrors the OOXML diff where a
# handful of shared-string timestamp values were bumped forward by ~45s.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview.Range("G2").Value = "2016-08-21 03:08:13"
$wsOverview.Range("G2").NumberFormat = $dateFormat

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file.
$wsZhCn.Range("H2").Value = "2016-08-21 03:08:08"
$wsZhCn.Range("H2").NumberFormat = $dateFormat
$wsZhCn.Range("K2").Value = "2016-08-21 03:08:26"
$wsZhCn.Range("K2").NumberFormat = $dateFormat

# de-de sheet: "Correspond Handoff Datetime" (shared with the Overview value)
# and "Correspond Handback DateTime" for the first file.
$wsDeDe.Range("H2").Value = "2016-08-21 03:08:13"
$wsDeDe.Range("H2").NumberFormat = $dateFormat
$wsDeDe.Range("K2").Value = "2016-08-21 03:08:33"
$wsDeDe.Range("K2").NumberFormat = $dateFormat
